$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 29303
$ws.Range("E2").Value = 2802
$ws.Range("F2").Value = 2802
$ws.Range("G2").Value = 1169
$ws.Range("H2").Value = 560
$ws.Range("I2").Value = 593
$ws.Range("J2").Value = -34
$ws.Range("K2").Value = 62563
$ws.Range("L2").Value = 32045
$ws.Range("M2").Value = 30519
$ws.Range("N2").Value = 26529
$ws.Range("O2").Value = 3989
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 2497
$ws.Range("R2").Value = -709
$ws.Range("S2").Value = -1992
$ws.Range("T2").Value = 528
$ws.Range("U2").Value = 1969
$ws.Range("V2").Value = 17242
$ws.Range("W2").Value = 9.56
$ws.Range("X2").Value = 1.91
$ws.Range("AA2").Value = 105
$ws.Range("AB2").Value = 5600786.91
$ws.Range("AC2").Value = 743
$ws.Range("AE2").Value = 30319
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 87500000
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AH2").ClearContents()

# Row 3
$ws.Range("D3").Value = 40408
$ws.Range("E3").Value = 3856
$ws.Range("F3").Value = 3856
$ws.Range("G3").Value = 2925
$ws.Range("H3").Value = 1481
$ws.Range("I3").Value = 1368
$ws.Range("J3").Value = 113
$ws.Range("K3").Value = 63649
$ws.Range("L3").Value = 32395
$ws.Range("M3").Value = 31255
$ws.Range("N3").Value = 27428
$ws.Range("O3").Value = 3827
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 3420
$ws.Range("R3").Value = -632
$ws.Range("S3").Value = -120
$ws.Range("T3").Value = 587
$ws.Range("U3").Value = 2833
$ws.Range("V3").Value = 18544
$ws.Range("W3").Value = 9.539999999999999
$ws.Range("X3").Value = 3.67
$ws.Range("Y3").Value = 5.07
$ws.Range("Z3").Value = 2.35
$ws.Range("AA3").Value = 103.65
$ws.Range("AB3").Value = 5714746.35
$ws.Range("AC3").Value = 1540
$ws.Range("AE3").Value = 31085
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 66613750
$ws.Range("AD3").ClearContents()
$ws.Range("AH3").ClearContents()

# Row 4
$ws.Range("D4").Value = 32873
$ws.Range("E4").Value = 3980
$ws.Range("F4").Value = 4140
$ws.Range("G4").Value = 2864
$ws.Range("H4").Value = 1803
$ws.Range("I4").Value = 1689
$ws.Range("J4").Value = 114
$ws.Range("K4").Value = 63261
$ws.Range("L4").Value = 29742
$ws.Range("M4").Value = 33518
$ws.Range("N4").Value = 33518
$ws.Range("P4").Value = 521
$ws.Range("Q4").Value = 3269
$ws.Range("R4").Value = -901
$ws.Range("S4").Value = -3028
$ws.Range("T4").Value = 524
$ws.Range("U4").Value = 2745
$ws.Range("V4").Value = 16121
$ws.Range("W4").Value = 12.11
$ws.Range("X4").Value = 5.49
$ws.Range("Y4").Value = 5.54
$ws.Range("Z4").Value = 2.84
$ws.Range("AA4").Value = 88.73
$ws.Range("AB4").Value = 7721.88
$ws.Range("AC4").Value = 1779
$ws.Range("AD4").Value = 20.15
$ws.Range("AE4").Value = 33435
$ws.Range("AF4").Value = 1.07
$ws.Range("AG4").Value = 700
$ws.Range("AH4").Value = 1.95
$ws.Range("AI4").Value = 41.54
$ws.Range("AJ4").Value = 100249166
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 33892
$ws.Range("E5").Value = 3945
$ws.Range("F5").Value = 3945
$ws.Range("G5").Value = 3148
$ws.Range("H5").Value = 2738
$ws.Range("I5").Value = 2738
$ws.Range("K5").Value = 61838
$ws.Range("L5").Value = 27563
$ws.Range("M5").Value = 34275
$ws.Range("N5").Value = 34275
$ws.Range("P5").Value = 462
$ws.Range("Q5").Value = 3766
$ws.Range("R5").Value = -1147
$ws.Range("S5").Value = -1790
$ws.Range("T5").Value = 469
$ws.Range("U5").Value = 3298
$ws.Range("V5").Value = 13385
$ws.Range("W5").Value = 11.64
$ws.Range("X5").Value = 8.08
$ws.Range("Y5").Value = 8.08
$ws.Range("Z5").Value = 4.38
$ws.Range("AA5").Value = 80.42
$ws.Range("AB5").Value = 8131.3
$ws.Range("AC5").Value = 2731
$ws.Range("AD5").Value = 13.11
$ws.Range("AE5").Value = 34190
$ws.Range("AF5").Value = 1.05
$ws.Range("AG5").Value = 800
$ws.Range("AH5").Value = 2.23
$ws.Range("AI5").Value = 29.3
$ws.Range("AJ5").Value = 100249166
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 39708
$ws.Range("E6").Value = 4590
$ws.Range("F6").Value = 4590
$ws.Range("G6").Value = 4037
$ws.Range("H6").Value = 2645
$ws.Range("I6").Value = 2645
$ws.Range("K6").Value = 64463
$ws.Range("L6").Value = 27067
$ws.Range("M6").Value = 37396
$ws.Range("N6").Value = 37396
$ws.Range("P6").Value = 482
$ws.Range("Q6").Value = 6166
$ws.Range("R6").Value = -679
$ws.Range("S6").Value = -4219
$ws.Range("T6").Value = 753
$ws.Range("U6").Value = 5413
$ws.Range("V6").Value = 10992
$ws.Range("W6").Value = 11.56
$ws.Range("X6").Value = 6.66
$ws.Range("Y6").Value = 7.38
$ws.Range("Z6").Value = 4.19
$ws.Range("AA6").Value = 72.38
$ws.Range("AB6").Value = 8538.209999999999
$ws.Range("AC6").Value = 2638
$ws.Range("AD6").Value = 11.94
$ws.Range("AE6").Value = 37303
$ws.Range("AF6").Value = 0.84
$ws.Range("AG6").Value = 900
$ws.Range("AH6").Value = 2.86
$ws.Range("AI6").Value = 34.11
$ws.Range("AJ6").Value = 100249166

# Row 7
$ws.Range("D7").Value = 44532
$ws.Range("E7").Value = 4818
$ws.Range("G7").Value = 4159
$ws.Range("H7").Value = 2933
$ws.Range("I7").Value = 2918
$ws.Range("K7").Value = 67065
$ws.Range("L7").Value = 27611
$ws.Range("M7").Value = 39456
$ws.Range("N7").Value = 39456
$ws.Range("P7").Value = 492
$ws.Range("Q7").Value = 3984
$ws.Range("R7").Value = -2238
$ws.Range("S7").Value = -3054
$ws.Range("T7").Value = 1096
$ws.Range("U7").Value = 2825
$ws.Range("W7").Value = 10.82
$ws.Range("X7").Value = 6.59
$ws.Range("Y7").Value = 7.59
$ws.Range("Z7").Value = 4.46
$ws.Range("AA7").Value = 69.98
$ws.Range("AC7").Value = 2910
$ws.Range("AD7").Value = 10.58
$ws.Range("AE7").Value = 39358
$ws.Range("AF7").Value = 0.78
$ws.Range("AG7").Value = 1120
$ws.Range("AH7").Value = 3.64
$ws.Range("AI7").Value = 38.48

# Row 8
$ws.Range("D8").Value = 45800
$ws.Range("E8").Value = 4979
$ws.Range("G8").Value = 4474
$ws.Range("H8").Value = 3266
$ws.Range("I8").Value = 3186
$ws.Range("K8").Value = 68145
$ws.Range("L8").Value = 26992
$ws.Range("M8").Value = 41153
$ws.Range("N8").Value = 41153
$ws.Range("P8").Value = 492
$ws.Range("Q8").Value = 4360
$ws.Range("R8").Value = -1146
$ws.Range("S8").Value = -2228
$ws.Range("T8").Value = 1001
$ws.Range("U8").Value = 3840
$ws.Range("W8").Value = 10.87
$ws.Range("X8").Value = 7.13
$ws.Range("Y8").Value = 7.9
$ws.Range("Z8").Value = 4.83
$ws.Range("AA8").Value = 65.59
$ws.Range("AC8").Value = 3178
$ws.Range("AD8").Value = 9.69
$ws.Range("AE8").Value = 41051
$ws.Range("AF8").Value = 0.75
$ws.Range("AG8").Value = 1196
$ws.Range("AH8").Value = 3.88
$ws.Range("AI8").Value = 37.62

# Row 9
$ws.Range("D9").Value = 47413
$ws.Range("E9").Value = 5197
$ws.Range("G9").Value = 4797
$ws.Range("H9").Value = 3414
$ws.Range("I9").Value = 3414
$ws.Range("K9").Value = 69916
$ws.Range("L9").Value = 26970
$ws.Range("M9").Value = 42946
$ws.Range("N9").Value = 42946
$ws.Range("P9").Value = 492
$ws.Range("Q9").Value = 5065
$ws.Range("R9").Value = -1270
$ws.Range("S9").Value = -2158
$ws.Range("T9").Value = 1014
$ws.Range("U9").Value = 4093
$ws.Range("W9").Value = 10.96
$ws.Range("X9").Value = 7.2
$ws.Range("Y9").Value = 8.119999999999999
$ws.Range("Z9").Value = 4.95
$ws.Range("AA9").Value = 62.8
$ws.Range("AC9").Value = 3405
$ws.Range("AD9").Value = 9.050000000000001
$ws.Range("AE9").Value = 42839
$ws.Range("AF9").Value = 0.72
$ws.Range("AG9").Value = 1236
$ws.Range("AH9").Value = 4.01
$ws.Range("AI9").Value = 36.3
